# Update the "想去人数" (interested-people count) figures in column F
# across the relevant worksheets, as produced by the latest scheduled
# gh-pages data refresh (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value  = 2826
$ws1.Range("F8").Value  = 1679
$ws1.Range("F9").Value  = 1821
$ws1.Range("F12").Value = 746
$ws1.Range("F13").Value = 897
$ws1.Range("F14").Value = 168
$ws1.Range("F15").Value = 377
$ws1.Range("F16").Value = 1115
$ws1.Range("F20").Value = 6632
$ws1.Range("F22").Value = 1535
$ws1.Range("F23").Value = 160
$ws1.Range("F26").Value = 309
$ws1.Range("F27").Value = 270
$ws1.Range("F30").Value = 911
$ws1.Range("F32").Value = 92
$ws1.Range("F34").Value = 707
$ws1.Range("F35").Value = 1468
$ws1.Range("F38").Value = 219
$ws1.Range("F39").Value = 16
$ws1.Range("F41").Value = 205
$ws1.Range("F42").Value = 159

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 12

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value  = 12
$ws4.Range("F10").Value = 2826
$ws4.Range("F11").Value = 1679
$ws4.Range("F12").Value = 1821
$ws4.Range("F15").Value = 746
$ws4.Range("F17").Value = 897
$ws4.Range("F18").Value = 168
$ws4.Range("F19").Value = 377
$ws4.Range("F20").Value = 1115
$ws4.Range("F23").Value = 6632
$ws4.Range("F25").Value = 1535
$ws4.Range("F27").Value = 160
$ws4.Range("F30").Value = 309
$ws4.Range("F31").Value = 270
$ws4.Range("F34").Value = 911
$ws4.Range("F36").Value = 92
$ws4.Range("F38").Value = 707
$ws4.Range("F39").Value = 1468
$ws4.Range("F42").Value = 219
$ws4.Range("F43").Value = 16
$ws4.Range("F45").Value = 205
$ws4.Range("F49").Value = 159
